# Auto-generated script applying scheduled-runner market-price updates
# to the Leve profit tables across all eight crafting-job worksheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 0
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H17").Value = 1254.4138
$ws.Range("J17").Value = 1262.0769
$ws.Range("L17").Value = 3786.2307
$ws.Range("N17").Value = -4122.2307
$ws.Range("H29").Value = 699
$ws.Range("I29").Value = 699
$ws.Range("K29").Value = 2097
$ws.Range("M29").Value = -1816
$ws.Range("H43").Value = 6334.5713
$ws.Range("J43").Value = 8633.333000000001
$ws.Range("L43").Value = 8633.333000000001
$ws.Range("N43").Value = -8771.333000000001
$ws.Range("H58").Value = 2228.0667
$ws.Range("J58").Value = 3999.8572
$ws.Range("L58").Value = 11999.5716
$ws.Range("N58").Value = -12299.5716
$ws.Range("H132").Value = 7161.1816
$ws.Range("I132").Value = 8497.839
$ws.Range("J132").Value = 3973.7693
$ws.Range("K132").Value = 25493.517
$ws.Range("L132").Value = 11921.3079
$ws.Range("M132").Value = -22963.517
$ws.Range("N132").Value = -16981.3079
$ws.Range("H137").Value = 7942.951
$ws.Range("I137").Value = 9597.931
$ws.Range("K137").Value = 28793.793
$ws.Range("M137").Value = -26243.793

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8311.333000000001
$ws.Range("I32").Value = 8896.666999999999
$ws.Range("K32").Value = 8896.666999999999
$ws.Range("M32").Value = -8609.666999999999
$ws.Range("H45").Value = 7720.091
$ws.Range("I45").Value = 6552.625
$ws.Range("K45").Value = 6552.625
$ws.Range("M45").Value = -6175.625
$ws.Range("H61").Value = 3491.971
$ws.Range("I61").Value = 3369.4822
$ws.Range("K61").Value = 3369.4822
$ws.Range("M61").Value = -3157.4822
$ws.Range("H74").Value = 3118.2
$ws.Range("I74").Value = 1984.6666
$ws.Range("K74").Value = 1984.6666
$ws.Range("M74").Value = -1110.6666
$ws.Range("H77").Value = 3118.2
$ws.Range("I77").Value = 1984.6666
$ws.Range("K77").Value = 9923.333000000001
$ws.Range("M77").Value = -5555.333000000001
$ws.Range("H132").Value = 3535.797
$ws.Range("I132").Value = 3179.9355
$ws.Range("K132").Value = 9539.806500000001
$ws.Range("M132").Value = -7009.806500000001
$ws.Range("H136").Value = 3491.971
$ws.Range("I136").Value = 3369.4822
$ws.Range("K136").Value = 10108.4466
$ws.Range("M136").Value = -7558.446599999999

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 13413.866
$ws.Range("I94").Value = 15062.409
$ws.Range("K94").Value = 15062.409
$ws.Range("M94").Value = -14611.409
$ws.Range("H99").Value = 58931.285
$ws.Range("J99").Value = 15599.75
$ws.Range("L99").Value = 15599.75
$ws.Range("N99").Value = -18595.75
$ws.Range("H107").Value = 2539.9583
$ws.Range("I107").Value = 2498.6365
$ws.Range("J107").Value = 2994.5
$ws.Range("K107").Value = 2498.6365
$ws.Range("L107").Value = 2994.5
$ws.Range("M107").Value = -578.6365000000001
$ws.Range("N107").Value = -6834.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5027.7617
$ws.Range("I31").Value = 4616.706
$ws.Range("J31").Value = 6774.75
$ws.Range("K31").Value = 4616.706
$ws.Range("L31").Value = 6774.75
$ws.Range("M31").Value = -4321.706
$ws.Range("N31").Value = -7364.75
$ws.Range("H34").Value = 5027.7617
$ws.Range("I34").Value = 4616.706
$ws.Range("J34").Value = 6774.75
$ws.Range("K34").Value = 4616.706
$ws.Range("L34").Value = 6774.75
$ws.Range("M34").Value = -4414.706
$ws.Range("N34").Value = -7178.75
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H58").Value = 1358.6571
$ws.Range("I58").Value = 1147.0968
$ws.Range("J58").Value = 2998.25
$ws.Range("K58").Value = 1147.0968
$ws.Range("L58").Value = 2998.25
$ws.Range("M58").Value = -944.0968
$ws.Range("N58").Value = -3404.25
$ws.Range("H132").Value = 16506.654
$ws.Range("I132").Value = 752.5833
$ws.Range("J132").Value = 205555.5
$ws.Range("K132").Value = 2257.7499
$ws.Range("L132").Value = 616666.5
$ws.Range("M132").Value = 272.2501000000002
$ws.Range("N132").Value = -621726.5
$ws.Range("H134").Value = 3715.7778
$ws.Range("I134").Value = 3052
$ws.Range("K134").Value = 9156
$ws.Range("M134").Value = -6621
$ws.Range("H136").Value = 1358.6571
$ws.Range("I136").Value = 1147.0968
$ws.Range("J136").Value = 2998.25
$ws.Range("K136").Value = 3441.2904
$ws.Range("L136").Value = 8994.75
$ws.Range("M136").Value = -891.2903999999999
$ws.Range("N136").Value = -14094.75

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 222
$ws.Range("I14").Value = 222
$ws.Range("K14").Value = 666
$ws.Range("M14").Value = -493
$ws.Range("H80").Value = 58108.668
$ws.Range("I80").Value = 2499.6
$ws.Range("K80").Value = 7498.799999999999
$ws.Range("M80").Value = -6562.799999999999
$ws.Range("H83").Value = 58108.668
$ws.Range("I83").Value = 2499.6
$ws.Range("K83").Value = 22496.4
$ws.Range("M83").Value = -17816.4
$ws.Range("H131").Value = 3341.8333
$ws.Range("I131").Value = 4529.7896
$ws.Range("J131").Value = 2014.1177
$ws.Range("K131").Value = 13589.3688
$ws.Range("L131").Value = 6042.3531
$ws.Range("M131").Value = -8549.3688
$ws.Range("N131").Value = -16122.3531
$ws.Range("H132").Value = 57175.668
$ws.Range("J132").Value = 73019.42999999999
$ws.Range("L132").Value = 657174.8699999999
$ws.Range("N132").Value = -662234.8699999999
$ws.Range("H136").Value = 2161.9285
$ws.Range("I136").Value = 741.4
$ws.Range("K136").Value = 2224.2
$ws.Range("M136").Value = 2875.8
$ws.Range("H141").Value = 3180.75
$ws.Range("I141").Value = 3157.6667
$ws.Range("J141").Value = 3250
$ws.Range("K141").Value = 9473.000100000001
$ws.Range("L141").Value = 9750
$ws.Range("M141").Value = -4293.000100000001
$ws.Range("N141").Value = -20110

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 15012.5
$ws.Range("J54").Value = 15012.5
$ws.Range("L54").Value = 15012.5
$ws.Range("N54").Value = -15792.5
$ws.Range("H80").Value = 12178.571
$ws.Range("I80").Value = 14422.889
$ws.Range("K80").Value = 14422.889
$ws.Range("M80").Value = -13424.889
$ws.Range("H83").Value = 12178.571
$ws.Range("I83").Value = 14422.889
$ws.Range("K83").Value = 72114.44499999999
$ws.Range("M83").Value = -67122.44499999999
$ws.Range("H97").Value = 3352.9019
$ws.Range("I97").Value = 3488.3333
$ws.Range("K97").Value = 3488.3333
$ws.Range("M97").Value = -2992.3333
$ws.Range("H113").Value = 55499.5
$ws.Range("I113").Value = 55499.5
$ws.Range("K113").Value = 55499.5
$ws.Range("M113").Value = -53329.5
$ws.Range("H123").Value = 37965.832
$ws.Range("J123").Value = 37965.832
$ws.Range("L123").Value = 37965.832
$ws.Range("N123").Value = -42865.832
$ws.Range("H132").Value = 5501.606
$ws.Range("I132").Value = 5482.727
$ws.Range("J132").Value = 5539.364
$ws.Range("K132").Value = 16448.181
$ws.Range("L132").Value = 16618.092
$ws.Range("M132").Value = -13918.181
$ws.Range("N132").Value = -21678.092

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H12").Value = 2500
$ws.Range("J12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("N12").Value = -1340
$ws.Range("H17").Value = 256.66666
$ws.Range("J17").Value = 500
$ws.Range("L17").Value = 500
$ws.Range("N17").Value = -840
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H55").Value = 1530.6154
$ws.Range("I55").Value = 378.42856
$ws.Range("K55").Value = 378.42856
$ws.Range("M55").Value = -205.42856
$ws.Range("H122").Value = 7126.3184
$ws.Range("I122").Value = 7419.2
$ws.Range("K122").Value = 22257.6
$ws.Range("M122").Value = -19807.6
$ws.Range("H132").Value = 307245.75
$ws.Range("I132").Value = 453313.56
$ws.Range("J132").Value = 5980.8125
$ws.Range("K132").Value = 1359940.68
$ws.Range("L132").Value = 17942.4375
$ws.Range("M132").Value = -1357410.68
$ws.Range("N132").Value = -23002.4375
$ws.Range("H136").Value = 4752
$ws.Range("I136").Value = 2005.5
$ws.Range("K136").Value = 6016.5
$ws.Range("M136").Value = -3466.5

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 30849.75
$ws.Range("I107").Value = 6381.5454
$ws.Range("K107").Value = 19144.6362
$ws.Range("M107").Value = -17224.6362
$ws.Range("H132").Value = 11107.6
$ws.Range("I132").Value = 12203.593
$ws.Range("J132").Value = 5727.273
$ws.Range("K132").Value = 36610.779
$ws.Range("L132").Value = 17181.819
$ws.Range("M132").Value = -34080.779
$ws.Range("N132").Value = -22241.819
$ws.Range("H136").Value = 315258
$ws.Range("I136").Value = 358701
$ws.Range("J136").Value = 3916.5
$ws.Range("K136").Value = 1076103
$ws.Range("L136").Value = 11749.5
$ws.Range("M136").Value = -1073553
$ws.Range("N136").Value = -16849.5

